$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.718.41"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "1.859.01"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("E4").Value = "  -0.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.670"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.26%  "
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.59"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.338"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0727"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "12.76"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "2.133.51"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.703"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "1.873.31"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.78"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").Value = "34.704.78"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("D19").Value = "0.0₃0805"
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.13"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -13.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.70"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.26"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.92"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("E29").Value = "  -5.85%  "
$ws.Range("D30").Value = "4.128.45"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0567"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.08"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.819"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -10.89%  "
$ws.Range("E37").Value = "  -20.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.92"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "96.89"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.76"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0661"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0208"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("E43").Value = "  -5.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0832"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +12.67%  "
$ws.Range("D45").Value = "1.275.65"
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.04%  "
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.74"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "41.99"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.02%  "
